$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.140.49"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.048.78"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'247.87"
$ws.Range("E5").Value = "  -2.06%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.69%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'57.14"
$ws.Range("E8").Value = "  -3.59%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.45%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0774"
$ws.Range("E10").Value = "  -3.17%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.108"
$ws.Range("E11").Value = "  +0.06%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "'15.77"
$ws.Range("E12").Value = "  -3.03%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.866"
$ws.Range("E13").Value = "  +5.42%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.345.34"
$ws.Range("E14").Value = "  -1.33%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +2.57%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.048.36"
$ws.Range("E16").Value = "  -1.24%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'17.87"
$ws.Range("E17").Value = "  +13.88%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.108.60"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "'74.54"
$ws.Range("E19").Value = "  -0.21%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -4.54%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'5.36"
$ws.Range("E21").Value = "  -2.06%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'236.60"
$ws.Range("E22").Value = "  -1.47%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.06%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +2.66%  "

# Row 25 - was PancakeSwap, now Cosmos
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.44"
$ws.Range("E25").Value = "  +1.01%  "

# Row 26 - was Cosmos, now PancakeSwap
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.17"
$ws.Range("E26").Value = "  -4.34%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'169.00"

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'20.03"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -1.59%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'4.81"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.79%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "'0.0616"
$ws.Range("E32").Value = "  -2.92%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.47%  "

# Row 34 - Kaspa
$ws.Range("D34").Value = "'0.0889"
$ws.Range("E34").Value = "  -1.79%  "

# Row 35 - BinanceUSD
$ws.Range("E35").Value = "  -0.01%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -2.00%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  -0.02%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  -2.15%  "

# Row 39 - HuobiToken
$ws.Range("E39").Value = "  +11.94%  "

# Row 40 - THORChain
$ws.Range("D40").Value = "'5.24"
$ws.Range("E40").Value = "  +16.31%  "

# Row 41 - Cronos
$ws.Range("D41").Value = "'0.0989"
$ws.Range("E41").Value = "  -14.70%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "'0.0221"
$ws.Range("E42").Value = "  -2.50%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "'17.14"
$ws.Range("E43").Value = "  -3.50%  "

# Row 44 - ARBITRUM
$ws.Range("E44").Value = "  -2.46%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'95.57"
$ws.Range("E45").Value = "  -3.39%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  -2.42%  "

# Row 47 - Maker
$ws.Range("D47").Value = "1.265.62"
$ws.Range("E47").Value = "  -2.84%  "

# Row 48 - MXToken
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -2.81%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -2.16%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.227.21"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "'43.53"
$ws.Range("E51").Value = "  -1.19%  "
